$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = -0.0314
$ws.Range("E5").Value = 0.4518
$ws.Range("H5").Value = -0.0846
$ws.Range("I5").Value = -0.0524
$ws.Range("L5").Value = -0.1381

# Row 6
$ws.Range("D6").Value = -0.0211
$ws.Range("E6").Value = -0.0166
$ws.Range("F6").Value = -0.0518
$ws.Range("G6").Value = -0.0493
$ws.Range("H6").Value = -0.0553
$ws.Range("I6").Value = -0.0403
$ws.Range("J6").Value = -0.0726
$ws.Range("K6").Value = -0.0683
$ws.Range("L6").Value = -0.0429
$ws.Range("M6").Value = -0.0348

# Row 7
$ws.Range("D7").Value = -0.0357
$ws.Range("E7").Value = -0.0373
$ws.Range("F7").Value = 0.1199
$ws.Range("G7").Value = 0.2176
$ws.Range("H7").Value = -0.0736
$ws.Range("I7").Value = -0.1049
$ws.Range("J7").Value = -0.1597
$ws.Range("K7").Value = -0.1803
$ws.Range("L7").Value = 0.044
$ws.Range("M7").Value = 0.1279

# Row 8
$ws.Range("D8").Value = -0.1929
$ws.Range("E8").Value = -0.1674
$ws.Range("F8").Value = -0.1227
$ws.Range("G8").Value = -0.396
$ws.Range("H8").Value = -0.3017
$ws.Range("I8").Value = -0.2013
$ws.Range("J8").Value = -0.0776
$ws.Range("K8").Value = -0.0344
$ws.Range("L8").Value = -0.0176
$ws.Range("M8").Value = 0.0018

# Row 14
$ws.Range("D14").Value = -0.339
$ws.Range("E14").Value = -1.1764
$ws.Range("F14").Value = -1.0187
$ws.Range("G14").Value = -1.5079
$ws.Range("H14").Value = -1.1877
$ws.Range("I14").Value = -0.5185
$ws.Range("K14").Value = -0.3408
$ws.Range("L14").Value = -0.2499

# Row 15
$ws.Range("D15").Value = -1.8671
$ws.Range("E15").Value = -2.922
$ws.Range("F15").Value = -2.1494
$ws.Range("G15").Value = -3.6555
$ws.Range("H15").Value = -2.1715
$ws.Range("I15").Value = -1.8151
$ws.Range("J15").Value = -1.7742
$ws.Range("K15").Value = -2.866
$ws.Range("L15").Value = -1.2664
$ws.Range("M15").Value = -0.4869

# Row 21
$ws.Range("D21").Value = 0.3351
$ws.Range("E21").Value = -0.0136
$ws.Range("F21").Value = 0.331
$ws.Range("G21").Value = 0.0586
$ws.Range("H21").Value = -0.0482
$ws.Range("J21").Value = -0.2683
$ws.Range("L21").Value = -0.0562

# Row 23
$ws.Range("D23").Value = -0.0246
$ws.Range("E23").Value = 0.1351
$ws.Range("F23").Value = -0.0017
$ws.Range("G23").Value = -0.0963
$ws.Range("H23").Value = 0.1275
$ws.Range("I23").Value = 0.124
$ws.Range("J23").Value = 0.2015
$ws.Range("K23").Value = 0.1896
$ws.Range("L23").Value = -0.0175
$ws.Range("M23").Value = 0.0083

# Row 24
$ws.Range("D24").Value = -0.2279
$ws.Range("E24").Value = -0.182
$ws.Range("F24").Value = -0.1537
$ws.Range("G24").Value = -0.211
$ws.Range("H24").Value = -0.1659
$ws.Range("I24").Value = -0.1707
$ws.Range("J24").Value = -0.1414
$ws.Range("K24").Value = -0.0676
$ws.Range("L24").Value = -0.0316
$ws.Range("M24").Value = -0.0192

# Row 28
$ws.Range("D28").Value = -0.2049
$ws.Range("E28").Value = -0.0568
$ws.Range("F28").Value = -0.0988
$ws.Range("G28").Value = -0.2135

# Row 32
$ws.Range("D32").Value = 0
$ws.Range("H32").Value = 0.0001
$ws.Range("I32").Value = 0.0001

# Row 33
$ws.Range("D33").Value = -0.0205
$ws.Range("E33").Value = -0.02
$ws.Range("F33").Value = 0.0006
$ws.Range("G33").Value = 0.0006
$ws.Range("H33").Value = 0.0006
$ws.Range("I33").Value = 0.0006
$ws.Range("J33").Value = 0.0006
$ws.Range("K33").Value = 0.0006
$ws.Range("L33").Value = 0.0006
$ws.Range("M33").Value = 0.0006

# Row 34
$ws.Range("D34").Value = 0.0042
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0.0621
$ws.Range("G34").Value = 0.0581
$ws.Range("H34").Value = 0.0438
$ws.Range("I34").Value = 0.0468
$ws.Range("J34").Value = -0.013
$ws.Range("K34").Value = -0.0097
$ws.Range("L34").Value = -0.0014
$ws.Range("M34").Value = 0.0004

# Row 35
$ws.Range("D35").Value = -0.1501
$ws.Range("E35").Value = -0.1215
$ws.Range("F35").Value = -0.1209
$ws.Range("G35").Value = -0.1195
$ws.Range("H35").Value = -0.1177
$ws.Range("I35").Value = -0.1162
$ws.Range("J35").Value = -0.0255
$ws.Range("K35").Value = -0.0018
$ws.Range("L35").Value = -0.0018
$ws.Range("M35").Value = -0.0013

# Row 41
$ws.Range("D41").Value = -0.0047
$ws.Range("E41").Value = 0.001
$ws.Range("F41").Value = 0.0004
$ws.Range("G41").Value = 0.0018
$ws.Range("H41").Value = 0.0026
$ws.Range("I41").Value = 0.001
$ws.Range("K41").Value = 0.0013
$ws.Range("L41").Value = 0.001

# Row 42
$ws.Range("D42").Value = -0.0898
$ws.Range("E42").Value = -0.0595
$ws.Range("F42").Value = -0.0443
$ws.Range("G42").Value = -0.1286
$ws.Range("H42").Value = -0.1382
$ws.Range("I42").Value = -0.1339
$ws.Range("J42").Value = -0.0396
$ws.Range("K42").Value = -0.0072
$ws.Range("L42").Value = -0.0023
$ws.Range("M42").Value = 0

# Row 48
$ws.Range("D48").Value = -0.0006
$ws.Range("E48").Value = -0.0006
$ws.Range("F48").Value = -0.0004
$ws.Range("G48").Value = -0.0003
$ws.Range("H48").Value = -0.0002
$ws.Range("J48").Value = 0.0001
$ws.Range("K48").Value = 0.0001
$ws.Range("L48").Value = 0.0001

# Row 50
$ws.Range("D50").Value = 0.0911
$ws.Range("E50").Value = 0.0879
$ws.Range("F50").Value = 0.0192
$ws.Range("G50").Value = -0.0664
$ws.Range("H50").Value = -0.0632
$ws.Range("I50").Value = -0.0619
$ws.Range("J50").Value = 0.0022
$ws.Range("K50").Value = 0.0032
$ws.Range("L50").Value = -0.0007
$ws.Range("M50").Value = 0.0002

# Row 51
$ws.Range("D51").Value = -0.0084
$ws.Range("E51").Value = -0.0049
$ws.Range("F51").Value = -0.0041
$ws.Range("G51").Value = -0.004
$ws.Range("H51").Value = -0.004
$ws.Range("I51").Value = -0.0039
$ws.Range("J51").Value = -0.0039
$ws.Range("K51").Value = -0.0007
$ws.Range("L51").Value = -0.0001

# Row 55
$ws.Range("D55").Value = 0.0004
$ws.Range("E55").Value = 0.0002
$ws.Range("F55").Value = 0.0002
$ws.Range("G55").Value = 0.0001
